# "full support to dropdown" - finish filling out the DropDown row in the
# UI-naming reference sheet: fix the control-type label casing, refresh the
# artlayer naming-prefix parameter list, add the sub-object naming
# convention, and mark the row as complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "下拉框@Dropdown"
$ws.Range("D9").Value = "b1_,b2_,b3,l1_,l2,m_"
$ws.Range("E9").Value = "vb_"
$ws.Range("F9").Value = "V"

# Zoom the view in and move the selection onto the row that was just
# completed.
$excel.ActiveWindow.Zoom = 115
$ws.Range("E9").Select()
